$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 with new address/coordinates
$ws.Range("A2").Value = "bhujangaraya sharma street,kavali"
$ws.Range("B2").Value = 14.913181
$ws.Range("C2").Value = 79.992981
$ws.Range("D2").Value = "bhujangaraya sharma street,kavali"

# Row 3: the address/coordinates that used to be in row 2
$ws.Range("A3").Value = "iit madras, guindy"
$ws.Range("B3").Value = 12.99466
$ws.Range("C3").Value = 80.23338
$ws.Range("D3").Value = "iit madras, guindy"

# Row 4: new address
$ws.Range("A4").Value = "kukatpally, hyderabad"
$ws.Range("B4").Value = 17.48846
$ws.Range("C4").Value = 78.40918000000001
$ws.Range("D4").Value = "kukatpally, hyderabad"

# Row 5: new address
$ws.Range("A5").Value = "sullurpeta, andhra pradesh"
$ws.Range("B5").Value = 13.70456
$ws.Range("C5").Value = 80.01612
$ws.Range("D5").Value = "sullurpeta, andhra pradesh"
